# The lecture-break rows (A308:A328, "Break"/"Break") originally had no
# Valence/Arousal ratings recorded. ReadData.py was updated to include all
# data files, so these rows now carry the same neutral Valence=3 / Arousal=3
# values used elsewhere for break periods.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D308:D328").Value = 3
$ws.Range("E308:E328").Value = 3

# Reflect where the editor's selection ended up after this pass.
$ws.Range("F326").Select()
